$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.516.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.16%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.883.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.74%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'244.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.31%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  +0.05%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.4692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.81%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.2892"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.53%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.06480"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.16%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'22.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.33%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.68%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.880.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.62%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "'Litecoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'95.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.26%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "'Polygon"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.7271"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.22%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'5.176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.74%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'282.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.84%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'30.503.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.10%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'12.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.19%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.000007470"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.31%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'2.130.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.72%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'5.242"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'6.230"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.10%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'163.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.01%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'9.063"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.40%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'18.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.62%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'1.887"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.91%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'1.332"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.90%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.09704"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.56%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.49%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'4.267"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.10%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'4.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.55%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.04862"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.73%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.57%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.6909"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.53%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.01887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.20%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'2.813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.36%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'75.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.39%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'6.160"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'2.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.25%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.4245"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.83%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.8232"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.23%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'101.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.05%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'9.533"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.96%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'35.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.35%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'6.953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.14%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'912.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.36%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.05753"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.89%  "
$ws.Range("E51").Style = "Normal"
